# Updated cryptos list on Thu Feb  1 15:23:48 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with latest scraped values.
# A leading apostrophe forces numeric-looking prices to stay plain text,
# matching the original inline-string cell type and avoiding float drift
# (e.g. "301.10" becoming 301.1 / 301.10000000000002).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.716.10"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "2.299.08"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'301.10"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").Value = "'96.82"
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("D7").Value = "'0.501"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "'33.69"
$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("D11").Value = "'0.0792"
$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").Value = "'48.81"
$ws.Range("E12").Value = "  -3.64%  "

$ws.Range("E13").Value = "  +2.42%  "

$ws.Range("D14").Value = "'16.50"
$ws.Range("E14").Value = "  +7.80%  "

$ws.Range("D15").Value = "'6.76"
$ws.Range("E15").Value = "  +1.90%  "

$ws.Range("D16").Value = "2.657.88"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").Value = "2.315.52"
$ws.Range("E17").Value = "  +1.23%  "

$ws.Range("D18").Value = "'0.797"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").Value = "42.652.42"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").Value = "'11.68"
$ws.Range("E20").Value = "  +1.49%  "

$ws.Range("D21").Value = "0.0₃0898"
$ws.Range("E21").Value = "  +0.60%  "

$ws.Range("D22").Value = "'6.02"
$ws.Range("E22").Value = "  +1.04%  "

$ws.Range("D23").Value = "'66.90"
$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("D24").Value = "'236.31"
$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("D25").Value = "'1.98"
$ws.Range("E25").Value = "  +2.41%  "

$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("D28").Value = "'24.05"
$ws.Range("E28").Value = "  -1.66%  "

$ws.Range("E29").Value = "  +0.56%  "

$ws.Range("D30").Value = "'34.56"
$ws.Range("E30").Value = "  +1.90%  "

$ws.Range("D31").Value = "'167.34"
$ws.Range("E31").Value = "  +2.64%  "

$ws.Range("D32").Value = "'9.17"
$ws.Range("E32").Value = "  +1.54%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("E34").Value = "  +8.07%  "

$ws.Range("D35").Value = "'4.98"
$ws.Range("E35").Value = "  +0.48%  "

$ws.Range("D36").Value = "'17.11"
$ws.Range("E36").Value = "  +3.89%  "

$ws.Range("E37").Value = "  -2.31%  "

$ws.Range("D38").Value = "'0.0697"
$ws.Range("E38").Value = "  +1.03%  "

$ws.Range("D39").Value = "'2.82"
$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("D40").Value = "'0.0999"
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").Value = "'1.75"
$ws.Range("E41").Value = "  -1.18%  "

$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").Value = "'2.37"
$ws.Range("E43").Value = "  -5.68%  "

$ws.Range("D44").Value = "1.965.44"
$ws.Range("E44").Value = "  +0.73%  "

$ws.Range("E45").Value = "  +1.08%  "

$ws.Range("E46").Value = "  -1.52%  "

$ws.Range("D47").Value = "'9.77"
$ws.Range("E47").Value = "  -3.17%  "

$ws.Range("D48").Value = "'2.84"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("D49").Value = "2.523.46"
$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("D50").Value = "'52.87"
$ws.Range("E50").Value = "  -2.92%  "

$ws.Range("D51").Value = "'1.50"
$ws.Range("E51").Value = "  +1.80%  "
